$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.773.61"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.493.75"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "3.488.04"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.191"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.02%  "
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "4.055.93"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "610.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "3.496.74"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "69.829.00"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.865"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -21.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.71%  "
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "623.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.60%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0473"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.12%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("D43").Value = "3.317.84"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "0.0₃0720"
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.308"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  -0.01%  "
